$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (all values in
# this column are plain text such as "42.273.95" or "0.631", never real
# numbers), so force Text format before writing the updated values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.273.95'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '2.246.88'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '246.83'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').Value = '77.15'
$ws.Range('E7').Value = '  +6.57%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.630'
$ws.Range('E9').Value = '  -2.40%  '
$ws.Range('D10').Value = '41.69'
$ws.Range('E10').Value = '  +6.90%  '
$ws.Range('D11').Value = '0.0945'
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('D14').Value = '2.580.40'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').Value = '14.95'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = '0.865'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '2.243.14'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').Value = '42.035.61'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('D19').Value = '0.0₃0981'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').Value = '71.58'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').Value = '231.85'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '2.18'
$ws.Range('E23').Value = '  -3.66%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  -5.30%  '
$ws.Range('D26').Value = '11.21'
$ws.Range('E26').Value = '  -3.45%  '
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').Value = '7.18'
$ws.Range('E28').Value = '  +11.41%  '
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').Value = '169.14'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').Value = '20.56'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('D32').Value = '0.0854'
$ws.Range('E32').Value = '  +5.59%  '
$ws.Range('D33').Value = '33.19'
$ws.Range('E33').Value = '  +6.82%  '
$ws.Range('D34').Value = '0.120'
$ws.Range('E34').Value = '  -5.85%  '
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('D36').Value = '4.62'
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('D37').Value = '4.90'
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('D38').Value = '0.0298'
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').Value = '13.40'
$ws.Range('E39').Value = '  -4.96%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '2.21'
$ws.Range('E40').Value = '  -4.70%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '5.92'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = '116.19'
$ws.Range('E42').Value = '  +19.49%  '
$ws.Range('D43').Value = '0.204'
$ws.Range('E43').Value = '  -4.82%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = '60.32'
$ws.Range('E44').Value = '  -2.54%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '8.81'
$ws.Range('E45').Value = '  -4.24%  '
$ws.Range('D46').Value = '0.101'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = '1.14'
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').Value = '4.27'
$ws.Range('E50').Value = '  -12.62%  '
$ws.Range('D51').Value = '2.29'
$ws.Range('E51').Value = '  +0.31%  '
